# Rename the default sheet so exported IC reference workbooks show a
# friendlier tab name instead of the generic "Sheet1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Query Results"
